$d = $word.ActiveDocument

# ==================================================================
# PHASE 1: append the new "Knärot" section (plain text only, no
#          character formatting yet) after the last paragraph
#          ("BILAGA 1 - Fridlysta arter"), right before the sectPr.
# ==================================================================
$insertionPara = $d.Paragraphs($d.Paragraphs.Count)
$insertionPara.Range.InsertParagraphAfter()

# --- new paragraph 1/13 (style: Heading1) ---
$newPara0 = $d.Paragraphs($d.Paragraphs.Count)
$newPara0.Style = "Heading1"
$r = $newPara0.Range
$r.Collapse(0)
$r.Text = "Knärot – ekologi samt krav på livsmiljön"

$tailPara = $d.Paragraphs($d.Paragraphs.Count)
$tailPara.Range.InsertParagraphAfter()

# --- new paragraph 2/13 (style: Normal) ---
$newPara1 = $d.Paragraphs($d.Paragraphs.Count)
$newPara1.Style = "Normal"
$r = $newPara1.Range
$r.Collapse(0)
$r.Text = "Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021)."

$tailPara = $d.Paragraphs($d.Paragraphs.Count)
$tailPara.Range.InsertParagraphAfter()

# --- new paragraph 3/13 (style: Normal) ---
$newPara2 = $d.Paragraphs($d.Paragraphs.Count)
$newPara2.Style = "Normal"
$r = $newPara2.Range
$r.Collapse(0)
$r.Text = "Samuel Johnsons doktorsavhandling “Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“ (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: “Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” Vidare “More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”"

$tailPara = $d.Paragraphs($d.Paragraphs.Count)
$tailPara.Range.InsertParagraphAfter()

# --- new paragraph 4/13 (style: Normal) ---
$newPara3 = $d.Paragraphs($d.Paragraphs.Count)
$newPara3.Style = "Normal"
$r = $newPara3.Range
$r.Collapse(0)
$r.Text = "Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: “In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”"

$tailPara = $d.Paragraphs($d.Paragraphs.Count)
$tailPara.Range.InsertParagraphAfter()

# --- new paragraph 5/13 (style: Normal) ---
$newPara4 = $d.Paragraphs($d.Paragraphs.Count)
$newPara4.Style = "Normal"
$r = $newPara4.Range
$r.Collapse(0)
$r.Text = "En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022)."

$tailPara = $d.Paragraphs($d.Paragraphs.Count)
$tailPara.Range.InsertParagraphAfter()

# --- new paragraph 6/13 (style: Normal) ---
$newPara5 = $d.Paragraphs($d.Paragraphs.Count)
$newPara5.Style = "Normal"
$r = $newPara5.Range
$r.Collapse(0)
$r.Text = "Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022)."

$tailPara = $d.Paragraphs($d.Paragraphs.Count)
$tailPara.Range.InsertParagraphAfter()

# --- new paragraph 7/13 (style: Heading2) ---
$newPara6 = $d.Paragraphs($d.Paragraphs.Count)
$newPara6.Style = "Heading2"
$r = $newPara6.Range
$r.Collapse(0)
$r.Text = "Referenser - knärot"

$tailPara = $d.Paragraphs($d.Paragraphs.Count)
$tailPara.Range.InsertParagraphAfter()

# --- new paragraph 8/13 (style: Normal) ---
$newPara7 = $d.Paragraphs($d.Paragraphs.Count)
$newPara7.Style = "Normal"
$r = $newPara7.Range
$r.Collapse(0)
$r.Text = "de Graaf M & Roberts M.R., 2009. Short-term response of the herbaceous layer within leave patches after harvest. Forest Ecology and Management 257, 1014-1025"

$tailPara = $d.Paragraphs($d.Paragraphs.Count)
$tailPara.Range.InsertParagraphAfter()

# --- new paragraph 9/13 (style: Normal) ---
$newPara8 = $d.Paragraphs($d.Paragraphs.Count)
$newPara8.Style = "Normal"
$r = $newPara8.Range
$r.Collapse(0)
$r.Text = "Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. Ecological Applications, 22, 2049-2064 "

$tailPara = $d.Paragraphs($d.Paragraphs.Count)
$tailPara.Range.InsertParagraphAfter()

# --- new paragraph 10/13 (style: Normal) ---
$newPara9 = $d.Paragraphs($d.Paragraphs.Count)
$newPara9.Style = "Normal"
$r = $newPara9.Range
$r.Collapse(0)
$r.Text = "Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. Interactive effects of drought and edge exposure on old-growth forest understory species. Landscape Ecology, 37, sid 1839-1853"

$tailPara = $d.Paragraphs($d.Paragraphs.Count)
$tailPara.Range.InsertParagraphAfter()

# --- new paragraph 11/13 (style: Normal) ---
$newPara10 = $d.Paragraphs($d.Paragraphs.Count)
$newPara10.Style = "Normal"
$r = $newPara10.Range
$r.Collapse(0)
$r.Text = "Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. Biological legacies buffer local species extinction after logging. Journal of Applied Ecology. 51, 53-62."

$tailPara = $d.Paragraphs($d.Paragraphs.Count)
$tailPara.Range.InsertParagraphAfter()

# --- new paragraph 12/13 (style: Normal) ---
$newPara11 = $d.Paragraphs($d.Paragraphs.Count)
$newPara11.Style = "Normal"
$r = $newPara11.Range
$r.Collapse(0)
$r.Text = "Skogsstyrelsen, 2022. Vägledning för hänsyn till knärot. https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/"

$tailPara = $d.Paragraphs($d.Paragraphs.Count)
$tailPara.Range.InsertParagraphAfter()

# --- new paragraph 13/13 (style: Normal) ---
$newPara12 = $d.Paragraphs($d.Paragraphs.Count)
$newPara12.Style = "Normal"
$r = $newPara12.Range
$r.Collapse(0)
$r.Text = "SLU Artdatabanken, 2021. Artfaktablad. Naturvård – artfakta. SLU Artdatabanken, Uppsala "

# ==================================================================
# PHASE 2: apply run-level character formatting (italics) now that
#          all paragraph text exists. Doing this in a separate pass
#          avoids the "current typing format" bleeding into later
#          freshly-typed plain-text paragraphs.
# ==================================================================
# paragraph 3/13
$paraStart = $newPara2.Range.Start
$paraEnd = $newPara2.Range.End
$cursor = $paraStart
$sr = $d.Range($cursor, $paraEnd)
[void]$sr.Find.Execute("“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“")
$sr.Font.Italic = $true
$cursor = $sr.End
$sr = $d.Range($cursor, $paraEnd)
[void]$sr.Find.Execute("“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ")
$sr.Font.Italic = $true
$cursor = $sr.End
$sr = $d.Range($cursor, $paraEnd)
[void]$sr.Find.Execute("“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”")
$sr.Font.Italic = $true
$cursor = $sr.End

# paragraph 4/13
$paraStart = $newPara3.Range.Start
$paraEnd = $newPara3.Range.End
$cursor = $paraStart
$sr = $d.Range($cursor, $paraEnd)
[void]$sr.Find.Execute("“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”")
$sr.Font.Italic = $true
$cursor = $sr.End

# paragraph 8/13
$paraStart = $newPara7.Range.Start
$paraEnd = $newPara7.Range.End
$cursor = $paraStart
$sr = $d.Range($cursor, $paraEnd)
[void]$sr.Find.Execute("Short-term response of the herbaceous layer within leave patches after harvest. ")
$sr.Font.Italic = $true
$cursor = $sr.End

# paragraph 9/13
$paraStart = $newPara8.Range.Start
$paraEnd = $newPara8.Range.End
$cursor = $paraStart
$sr = $d.Range($cursor, $paraEnd)
[void]$sr.Find.Execute("Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ")
$sr.Font.Italic = $true
$cursor = $sr.End

# paragraph 10/13
$paraStart = $newPara9.Range.Start
$paraEnd = $newPara9.Range.End
$cursor = $paraStart
$sr = $d.Range($cursor, $paraEnd)
[void]$sr.Find.Execute("Interactive effects of drought and edge exposure on old-growth forest understory species. ")
$sr.Font.Italic = $true
$cursor = $sr.End

# paragraph 11/13
$paraStart = $newPara10.Range.Start
$paraEnd = $newPara10.Range.End
$cursor = $paraStart
$sr = $d.Range($cursor, $paraEnd)
[void]$sr.Find.Execute("Biological legacies buffer local species extinction after logging. ")
$sr.Font.Italic = $true
$cursor = $sr.End

# paragraph 12/13
$paraStart = $newPara11.Range.Start
$paraEnd = $newPara11.Range.End
$cursor = $paraStart
$sr = $d.Range($cursor, $paraEnd)
[void]$sr.Find.Execute("Vägledning för hänsyn till knärot. ")
$sr.Font.Italic = $true
$cursor = $sr.End

# paragraph 13/13
$paraStart = $newPara12.Range.Start
$paraEnd = $newPara12.Range.End
$cursor = $paraStart
$sr = $d.Range($cursor, $paraEnd)
[void]$sr.Find.Execute("Artfaktablad. Naturvård – artfakta. ")
$sr.Font.Italic = $true
$cursor = $sr.End

# ==================================================================
# Update the date in the "first page" header (header3.xml).
# ==================================================================
$sec = $d.Sections(1)
$firstPageHeader = $sec.Headers(2)
[void]$firstPageHeader.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2)
